$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data row for 2025-01-07 21:33:52 (resale numbers update).
# Columns A-D are text-like ("01" week code must not become the number 1,
# the date must not become a date serial) so force text format first, then
# restore the default "Normal" style afterward so no stray style index is
# left attached to the new cells (matches the plain, unstyled data rows
# already in the sheet).
$ws.Range("A13:D13").NumberFormat = "@"

$ws.Range("A13").Value = "2025-01-07"
$ws.Range("B13").Value = "21:33:52"
$ws.Range("C13").Value = "Tuesday"
$ws.Range("D13").Value = "01"

$ws.Range("A13:D13").Style = "Normal"

# Numeric columns
$ws.Range("E13").Value = 127364
$ws.Range("F13").Value = 143619
$ws.Range("G13").Value = 169094
$ws.Range("H13").Value = 158968
$ws.Range("I13").Value = -1
$ws.Range("J13").Value = 142281
$ws.Range("K13").Value = -1
$ws.Range("L13").Value = -1
$ws.Range("M13").Value = 192383
$ws.Range("N13").Value = 115065
$ws.Range("O13").Value = 45737
$ws.Range("P13").Value = 28399
$ws.Range("Q13").Value = 64460
$ws.Range("R13").Value = -1
$ws.Range("S13").Value = 47767
$ws.Range("T13").Value = -1
